$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Length -ge 2 -and $parts[0].ToLower() -eq "system") {
            $first = $parts[0]
            $lastIdx = $parts.Length - 1
            $last = $parts[$lastIdx]
            $parts[0] = $last
            $parts[$lastIdx] = $first
            $newVal = $parts -join ", "
            $cell.Value = $newVal
        }
    }
}
